# Update Pins: add two new rows to the "Dashboard" RPM sensor table on Tabelle1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 34: RPM Sensor hinten Links -> PC_4 (rechts)
$ws.Range("C34").Value = "RPM Sensor hinten Links"
$ws.Range("C34").Style = "Gut"
$ws.Range("D34").Value = "PC_4 (rechts)"

# Row 35: RPM Sensor hinten Rechts -> PC_5 (rechts)
$ws.Range("C35").Value = "RPM Sensor hinten Rechts"
$ws.Range("C35").Style = "Gut"
$ws.Range("D35").Value = "PC_5 (rechts)"

# Scroll back to the top of the sheet and move the selection to D36,
# matching the author's final cursor position after entering the new rows.
$ws.Application.Goto($ws.Range("A1")) | Out-Null
$ws.Range("D36").Select() | Out-Null
